$d = $word.ActiveDocument

# Locate the paragraph that ends the "Man ska kunna ändra ..." bullet item
# in the "Framtida funktioner" list (numId 5) so we can insert a new
# sibling list item right after it, matching its formatting.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Man ska kunna ändra inställningarna av sin editor.*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

# Insert a brand-new paragraph right after the target; Word automatically
# continues the same list (ListParagraph style + numId 5, ilvl 0) because
# the new paragraph inherits paragraph formatting from the one before it.
$newRange = $target.Range.InsertParagraphAfter()

# The freshly inserted paragraph is the one following the target paragraph.
$newPara = $target.Next()
$newPara.Range.Text = "Se till att programmeringsspråksmenyn visar vilket språk man valt."
